$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 15-74: unit price (E) becomes 39293, discount (F) cleared
for ($r = 15; $r -le 74; $r++) {
    $ws.Cells.Item($r, 5).Value = 39293
    $ws.Cells.Item($r, 6).ClearContents()
}

# Rows 75-104: unit price (E) becomes 45105, discount (F) cleared
for ($r = 75; $r -le 104; $r++) {
    $ws.Cells.Item($r, 5).Value = 45105
    $ws.Cells.Item($r, 6).ClearContents()
}

$wb.Application.Calculate()
